$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "Inter TV Rural"
$ws.Cells.Item($row, 3).Value = "Agricultura"
$ws.Cells.Item($row, 4).Value = "2025-03-16T00:00"
$ws.Cells.Item($row, 5).Value = "Positivo"
$ws.Cells.Item($row, 6).Value = "Produção de soja em Campos. Agricultores do Norte Fluminense investem no grão tipo exportação. Campos é o maior produtor de soja no estado do Rio. Última safra foi aproximadamente 3 mil toneladas. Cerca de 850 hectares. Proximidade com o Porto do Açu é um dos benefícios pelo baixo valor do frete. Em Santa Cruz, está a maior concentração de soja da região. Entrevista com produtor rural, José Geraldo Neto. Entrevista com secretário de Agricultura, Almy Júnior e com engenheiro agrônomo, Elias Deulefeu. "
